$wb = $excel.ActiveWorkbook

$accounts = $wb.Worksheets.Item("Accounts")
$accounts.Range("C2").Value = 2995
$accounts.Range("C4").Value = 1785

$sales = $wb.Worksheets.Item("Sales")
$sales.Range("B2").Value = 325
$sales.Range("B3").Value = 850
